$d = $word.ActiveDocument

# The document ends with a paragraph "Sanath, Vignesh,Theerthesh".
# Append a brand-new paragraph right after it (before the sectPr),
# carrying the same run/paragraph formatting (en-US language), with
# the text "Hello, How are you".
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Hello, How are you"

Write-Output "Appended new paragraph"
